$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---- Row 2 (ethyl benzene) ----
$ws.Range("B2").Value = 1.59474299178
Set-TextValue $ws.Range("C2") "1.6e+00"
$ws.Range("P2").Value = 527
$ws.Range("Q2").Value = 187
$ws.Range("R2").Value = 230
$ws.Range("S2").Value = 230
Set-TextValue $ws.Range("T2") "01009"
$ws.Range("U2").Value = "1009U01294"
$ws.Range("V2").Value = 34
$ws.Range("W2").Value = 3776135
$ws.Range("X2").Value = 34.12597
$ws.Range("Y2").Value = -86.842941
$ws.Range("Z2").Value = "C"
$ws.Range("AA2").Value = "Discrete"

# ---- Row 3 (toluene) ----
$ws.Range("B3").Value = 1.35402706849
Set-TextValue $ws.Range("C3") "1.4e+00"
$ws.Range("P3").Value = 527
$ws.Range("Q3").Value = 187
$ws.Range("R3").Value = 230
$ws.Range("S3").Value = 230
Set-TextValue $ws.Range("T3") "01009"
$ws.Range("U3").Value = "1009U01294"
$ws.Range("V3").Value = 34
$ws.Range("W3").Value = 3776135
$ws.Range("X3").Value = 34.12597
$ws.Range("Y3").Value = -86.842941
$ws.Range("Z3").Value = "C"
$ws.Range("AA3").Value = "Discrete"

# ---- Row 4 (xylenes (mixed)) ----
$ws.Range("B4").Value = 46.8192470795
Set-TextValue $ws.Range("C4") "4.7e+01"
$ws.Range("P4").Value = 527
$ws.Range("Q4").Value = 187
$ws.Range("R4").Value = 230
$ws.Range("S4").Value = 230
Set-TextValue $ws.Range("T4") "01009"
$ws.Range("U4").Value = "1009U01294"
$ws.Range("V4").Value = 34
$ws.Range("W4").Value = 3776135
$ws.Range("X4").Value = 34.12597
$ws.Range("Y4").Value = -86.842941
$ws.Range("Z4").Value = "C"
$ws.Range("AA4").Value = "Discrete"
